# Samplify HowTo.docx update:
#  - add a new "Beta key controls ..." heading paragraph after "Samplify"
#  - tweak the g-> cue-point line to mention "last SET cue point"
#    (and carry the _GoBack bookmark along with it)
#  - add "drag sample into DAW..." after the right-click line
#  - drop the "Tag Explorer" label paragraph
#  - add a trailing "right click for more options" line
#  - drop everything from the old "New Tags:" section onward

$d = $word.ActiveDocument

function Get-ParaIndexByText($doc, $text) {
    $idx = -1
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $t = $doc.Paragraphs.Item($i).Range.Text.TrimEnd([char]13)
        if ($t -eq $text) {
            $idx = $i
        }
    }
    return $idx
}

# --- Step 1: drop the old tail: everything from the blank paragraph right
#     after "click tags to search for samples containing tag" through the
#     end of the document ("Repair SampleLibrary", which also carries the
#     old _GoBack bookmark).
$clickTagsIdx = Get-ParaIndexByText $d "click tags to search for samples containing tag"
$tailStartIdx = $clickTagsIdx + 1
$tailStartPara = $d.Paragraphs.Item($tailStartIdx)
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$tailRange = $d.Range($tailStartPara.Range.Start, $lastPara.Range.End)
$tailRange.Delete()

# --- Step 2: drop the "Tag Explorer" label paragraph entirely.
$tagExplorerIdx = Get-ParaIndexByText $d "Tag Explorer"
$d.Paragraphs.Item($tagExplorerIdx).Range.Delete()

# --- Step 3: add "drag sample into DAW..." right after the "right click on
#     sample thumbnail..." paragraph.
$rightClickIdx = Get-ParaIndexByText $d "right click on sample thumbnail-> play sample starting at selected point in sample"
$d.Paragraphs.Item($rightClickIdx).Range.InsertParagraphAfter()
$dragDawPara = $d.Paragraphs.Item($rightClickIdx + 1)
$dragDawPara.Range.Text = "drag sample into DAW when a good sound was found"

# --- Step 4: add "right click for more options" as a new last paragraph,
#     right after "click tags to search for samples containing tag".
$clickTagsIdx = Get-ParaIndexByText $d "click tags to search for samples containing tag"
$d.Paragraphs.Item($clickTagsIdx).Range.InsertParagraphAfter()
$moreOptionsPara = $d.Paragraphs.Item($clickTagsIdx + 1)
$moreOptionsPara.Range.Text = "right click for more options"

# --- Step 5: rewrite the "g-> start sample at last cue point" paragraph as
#     three runs ("...at last" / " set" / " cue point") with the _GoBack
#     bookmark re-inserted between the 2nd and 3rd run. Use InsertXML on the
#     paragraph text (minus the trailing paragraph mark) so the existing
#     runs/paragraph properties outside the replaced span stay untouched.
$gIdx = Get-ParaIndexByText $d "g-> start sample at last cue point"
$gPara = $d.Paragraphs.Item($gIdx)
$gRange = $gPara.Range
$gRangeNoMark = $d.Range($gRange.Start, $gRange.End - 1)
$gXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>g-&gt; start sample at last</w:t></w:r><w:r><w:t xml:space="preserve"> set</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> cue point</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$gRangeNoMark.InsertXML($gXml)

# --- Step 6: insert the new "Beta key controls (Customization and
#     remapping in works):" paragraph right after "Samplify", split into
#     its three runs via InsertXML (use the *whole* new-paragraph range,
#     mark included, since the paragraph is brand new/empty).
$samplifyIdx = Get-ParaIndexByText $d "Samplify"
$d.Paragraphs.Item($samplifyIdx).Range.InsertParagraphAfter()
$betaPara = $d.Paragraphs.Item($samplifyIdx + 1)
$betaXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Beta key controls (Customization and </w:t></w:r><w:r><w:t>remapping</w:t></w:r><w:r><w:t xml:space="preserve"> in works):</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$betaPara.Range.InsertXML($betaXml)

Write-Output "done"
